$d = $word.ActiveDocument

# Locate the paragraph that starts the "Additional helpful graphs..." text.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Additional helpful graphs*") {
        $target = $p
        break
    }
}

if ($target -ne $null) {
    $paraStart = $target.Range.Start
    $paraEnd   = $target.Range.End            # includes the end-of-paragraph mark

    $keepPrefix = "Additional helpful graphs/ tables would be to look at the impacts of "
    $oldRun1    = "Additional helpful graphs/ tables would be to look at the impacts of spotlight and staff pick on the average donation. "

    $newTail = "the average goal on the campaign outcome. It was found that in all parent categories except publishing and technology the average goal was lowest for successful campaigns. This means that future campaigns should attempt to keep costs down and have the lowest goal possible for their project."

    # End of the original first run (start of the text that used to be the 2nd run).
    $run1End = $paraStart + $oldRun1.Length

    # Delete everything after the first run up to (but excluding) the paragraph mark
    # -- this removes the old 2nd/3rd/4th runs (incl. the "backers" proofErr-wrapped run).
    $restRange = $d.Range($run1End, $paraEnd - 1)
    $restRange.Delete()

    # Trim the tail of the first run down to "...impacts of " (still the very same run).
    $tailStart = $paraStart + $keepPrefix.Length
    $tailRange = $d.Range($tailStart, $run1End)
    $tailRange.Delete()

    # Insert the new second sentence right after the trimmed first run.
    $insertPoint = $d.Range($tailStart, $tailStart)
    $insertPoint.InsertAfter($newTail)

    # Force Word to keep the inserted text as its own run (rather than silently
    # re-merging it into the preceding, identically-formatted run) by toggling a
    # character property on just the new text.
    $newRunRange = $d.Range($tailStart, $target.Range.End - 1)
    $newRunRange.Font.Bold = $true
    $newRunRange.Font.Bold = $false
}
